$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)
$para = $sh.TextFrame.TextRange.Paragraphs(5)
$para.Runs(1).Text = "Payload length indicator (PLI) and header error control (HEC) are for tracing and debugging purposes interesting but not relevant for accounting. Therefore, excluded from IPFIX. "
